# Commit: "minor adjustments to emx"
# Adds freeze1-patch1, freeze1-patch3, freeze2-patch1 rows (and fixes/adds
# freeze3) to the "entities" sheet, together with new description text
# for all freeze* rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("entities")

# --- Insert two new rows right after the existing "freeze1" row (row 5) ---
$ws.Rows.Item(6).Insert()
$ws.Rows.Item(6).Insert()

# --- Insert one new row right after the existing "freeze2" row             ---
# (freeze2 was row 6 before inserting, now it is row 8 after the 2 inserts)
$ws.Rows.Item(9).Insert()

# Row layout is now:
#  1 header
#  2 phenopacket
#  3 ped
#  4 template
#  5 freeze1
#  6 <new: freeze1-patch1>
#  7 <new: freeze1-patch3>
#  8 freeze2
#  9 <new: freeze2-patch1>
# 10 freeze3

# --- freeze1 : add description ---
$ws.Cells.Item(5, 4).Value = "Original gvcf, bam, ped, phenopacket"

# --- freeze1-patch1 (new row 6) ---
$ws.Cells.Item(6, 1).Value = "rd3_portal_cluster"
$ws.Cells.Item(6, 2).Value = "freeze1-patch1"
$ws.Cells.Item(6, 3).Value = "Freeze 1 Patch 1 Files"
$ws.Cells.Item(6, 4).Value = "Updated ped and phenopacket files"
$ws.Cells.Item(6, 6).Value = "rd3_portal_cluster_template"

# --- freeze1-patch3 (new row 7) ---
$ws.Cells.Item(7, 1).Value = "rd3_portal_cluster"
$ws.Cells.Item(7, 2).Value = "freeze1-patch3"
$ws.Cells.Item(7, 3).Value = "Freeze 1 Patch 3 Files"
$ws.Cells.Item(7, 4).Value = "Updated ped and phenopacket files"
$ws.Cells.Item(7, 6).Value = "rd3_portal_cluster_template"

# --- freeze2 (row 8) : add description ---
$ws.Cells.Item(8, 4).Value = "Original gvcf, bam, ped, phenopacket"

# --- freeze2-patch1 (new row 9) ---
$ws.Cells.Item(9, 1).Value = "rd3_portal_cluster"
$ws.Cells.Item(9, 2).Value = "freeze2-patch1"
$ws.Cells.Item(9, 3).Value = "Freeze 2 Patch 1 Files"
$ws.Cells.Item(9, 4).Value = "Updated ped and phenopacket files"
$ws.Cells.Item(9, 6).Value = "rd3_portal_cluster_template"

# --- freeze3 (row 10) : fix label + add description ---
$ws.Cells.Item(10, 3).Value = "Freeze 3 Files"
$ws.Cells.Item(10, 4).Value = "Original gvcf, bam, ped, phenopacket"

$wb.Save()
